$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.246.14'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.859.57'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9994'
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7112'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.57'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9997'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07763'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3100'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.90'
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07811'
$ws.Range('E11').Value = '  -2.56%  '
$ws.Range('D12').Value = '1.862.51'
$ws.Range('E12').Value = '  -11.73%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.108'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.20'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6883'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.524'
$ws.Range('E16').Value = '  +3.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008435'
$ws.Range('E17').Value = '  +2.38%  '
$ws.Range('D18').Value = '29.240.69'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.46'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '2.107.32'
$ws.Range('E20').Value = '  -4.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.85'
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9994'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.527'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9995'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1542'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.46'
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.871'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.53'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.565'
$ws.Range('E29').Value = '  +4.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.254'
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.234'
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.197'
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05213'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7590'
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.846'
$ws.Range('E35').Value = '  -1.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.166'
$ws.Range('E36').Value = '  +0.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.708'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01861'
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('D39').Value = '1.226.56'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8971'
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.85'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9993'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.579'
$ws.Range('E44').Value = '  -10.08%  '
$ws.Range('D45').Value = '2.005.34'
$ws.Range('E45').Value = '  -2.53%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000124'
$ws.Range('E46').Value = '  -3.87%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5177'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '64.88'
$ws.Range('E48').Value = '  -9.86%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.525'
$ws.Range('E49').Value = '  +1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.750'
$ws.Range('E50').Value = '  -2.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.013'
$ws.Range('E51').Value = '  +0.58%  '
